$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 35; this shifts existing rows 35-137
# down to 36-138 (matching the diff, which re-bases every record after the
# old row 34 by one row and appends the former last row as the new row 138).
$ws.Rows(35).Insert()

# Populate the newly inserted row 35 with the new record from the diff.
$ws.Range("A35").Value = 5
$ws.Range("B35").Value = "Macroferia Regional de Talca"
$ws.Range("C35").Value = "Maule"
$ws.Range("D35").Value = 44414
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = 100112003
$ws.Range("G35").Value = "Ajo"
$ws.Range("H35").Value = "Chino"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 14000
$ws.Range("L35").Value = 14000
$ws.Range("M35").Value = 14000
$ws.Range("N35").Value = "$/caja 10 kilos"
$ws.Range("O35").Value = "China"
$ws.Range("P35").Value = 1400
$ws.Range("Q35").Value = 10
$ws.Range("R35").Value = "Hortaliza"
